$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (tab name)
$ws.Name = "Table_Boss"

# Update the selection to D22
$ws.Range("D22").Select() | Out-Null
